{"js": "const replacements = [\n  [\"120\u00d73=\", \"334\u00d72=\"],\n  [\"386\u00d73=\", \"947\u00d78=\"],\n  [\"382\u00d78=\", \"289\u00d73=\"],\n  [\"838\u00d74=\", \"206\u00d75=\"],\n  [\"476\u00d78=\", \"491\u00d79=\"],\n  [\"732\u00d76=\", \"314\u00d75=\"],\n  [\"530\u00d78=\", \"847\u00d73=\"],\n  [\"888\u00d79=\", \"128\u00d75=\"],\n  [\"209\u00d72=\", \"914\u00d75=\"],\n  [\"563\u00d75=\", \"125\u00d73=\"],\n  [\"856\u00d73=\", \"831\u00d76=\"],\n  [\"161\u00d79=\", \"650\u00d78=\"],\n  [\"415\u00d76=\", \"476\u00d74=\"],\n  [\"225\u00d74=\", \"847\u00d72=\"],\n  [\"843\u00d72=\", \"826\u00d74=\"],\n  [\"355\u00d72=\", \"658\u00d79=\"],\n  [\"440\u00d75=\", \"479\u00d74=\"],\n  [\"643\u00d73=\", \"899\u00d74=\"],\n  [\"434\u00d78=\", \"646\u00d79=\"],\n  [\"497\u00d72=\", \"885\u00d76=\"],\n  [\"937\u00d72=\", \"877\u00d75=\"],\n  [\"240\u00d75=\", \"985\u00d74=\"],\n  [\"340\u00d79=\", \"616\u00d73=\"],\n  [\"901\u00d74=\", \"627\u00d79=\"],\n  [\"832\u00d76=\", \"160\u00d73=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"120\u00d73=\", \"334\u00d72=\"),\n    @(\"386\u00d73=\", \"947\u00d78=\"),\n    @(\"382\u00d78=\", \"289\u00d73=\"),\n    @(\"838\u00d74=\", \"206\u00d75=\"),\n    @(\"476\u00d78=\", \"491\u00d79=\"),\n    @(\"732\u00d76=\", \"314\u00d75=\"),\n    @(\"530\u00d78=\", \"847\u00d73=\"),\n    @(\"888\u00d79=\", \"128\u00d75=\"),\n    @(\"209\u00d72=\", \"914\u00d75=\"),\n    @(\"563\u00d75=\", \"125\u00d73=\"),\n    @(\"856\u00d73=\", \"831\u00d76=\"),\n    @(\"161\u00d79=\", \"650\u00d78=\"),\n    @(\"415\u00d76=\", \"476\u00d74=\"),\n    @(\"225\u00d74=\", \"847\u00d72=\"),\n    @(\"843\u00d72=\", \"826\u00d74=\"),\n    @(\"355\u00d72=\", \"658\u00d79=\"),\n    @(\"440\u00d75=\", \"479\u00d74=\"),\n    @(\"643\u00d73=\", \"899\u00d74=\"),\n    @(\"434\u00d78=\", \"646\u00d79=\"),\n    @(\"497\u00d72=\", \"885\u00d76=\"),\n    @(\"937\u00d72=\", \"877\u00d75=\"),\n    @(\"240\u00d75=\", \"985\u00d74=\"),\n    @(\"340\u00d79=\", \"616\u00d73=\"),\n    @(\"901\u00d74=\", \"627\u00d79=\"),\n    @(\"832\u00d76=\", \"160\u00d73=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute([ref]$pair[0], $false, $false, $false, $false, $false, $true, 1, $false, [ref]$pair[1], 2)\n}"}
